$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet2: append the 2019 data rows (rows 620-654)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()

$rows2019 = @(
  @(1,  "State",           "Andhra Pradesh",    2019, 237567),
  @(2,  "State",           "Arunachal Pradesh", 2019, 2590),
  @(3,  "State",           "Assam",             2019, 123512),
  @(4,  "State",           "Bihar",             2019, 197935),
  @(5,  "State",           "Chhattisgarh",      2019, 61256),
  @(6,  "State",           "Goa",               2019, 2465),
  @(7,  "State",           "Gujarat",           2019, 139503),
  @(8,  "State",           "Haryana",           2019, 111323),
  @(9,  "State",           "Himachal Pradesh",  2019, 14480),
  @(10, "State",           "Jammu & Kashmir",   2019, 22404),
  @(11, "State",           "Jharkhand",         2019, 50048),
  @(12, "State",           "Karnataka",         2019, 120165),
  @(13, "State",           "Kerala",            2019, 175810),
  @(14, "State",           "Madhya Pradesh",    2019, 246470),
  @(15, "State",           "Maharashtra",       2019, 341084),
  @(16, "State",           "Manipur",           2019, 2830),
  @(17, "State",           "Meghalaya",         2019, 3125),
  @(18, "State",           "Mizoram",           2019, 2379),
  @(19, "State",           "Nagaland",          2019, 1117),
  @(20, "State",           "Odisha",            2019, 96033),
  @(21, "State",           "Punjab",            2019, 44697),
  @(22, "State",           "Rajasthan",         2019, 225306),
  @(23, "State",           "Sikkim",            2019, 632),
  @(24, "State",           "Tamil Nadu",        2019, 168116),
  @(25, "State",           "Tripura",           2019, 5336),
  @(26, "State",           "Uttar Pradesh",     2019, 353131),
  @(27, "State",           "Uttarakhand",       2019, 12081),
  @(28, "State",           "West Bengal",       2019, 157610),
  @(29, "Union Territory",  "A & N Islands",    2019, 564),
  @(30, "Union Territory",  "Chandigarh",       2019, 2819),
  @(31, "Union Territory",  "D & N Haveli",     2019, 226),
  @(32, "Union Territory",  "Daman & Diu",      2019, 322),
  @(33, "Union Territory",  "Delhi",            2019, 299475),
  @(34, "Union Territory",  "Lakshadweep",      2019, 123),
  @(35, "Union Territory",  "Puducherry",       2019, 3167)
)

$startRow = 620
for ($i = 0; $i -lt $rows2019.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows2019[$i]
    $ws2.Cells.Item($r, 1).Value = $vals[0]
    $ws2.Cells.Item($r, 2).Value = $vals[1]
    $ws2.Cells.Item($r, 3).Value = $vals[2]
    $ws2.Cells.Item($r, 4).Value = $vals[3]
    $ws2.Cells.Item($r, 5).Value = $vals[4]
}

$lastRow = $startRow + $rows2019.Count - 1
$ws2.Range("E" + $lastRow).Select()

# ---------------------------------------------------------------------------
# Sheet1: move the selection (no data change)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("G12").Select()

# Re-activate Sheet2 last so it remains the active tab, matching the target
$ws2.Activate()
$ws2.Range("E" + $lastRow).Select()
